$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Credit the sheet with a fresh batch of msisdn numbers (rows 2-6), keeping
# them as text (they're phone numbers, not numeric values) the same way a
# user would by prefixing the entry with an apostrophe.
$ws.Range("A2").Value = "'7075812222"
$ws.Range("A3").Value = "'9885861677"
$ws.Range("A4").Value = "'8341325077"
$ws.Range("A5").Value = "'9640494242"
$ws.Range("A6").Value = "'7014118238"

# Cut the rest of the old list away - only six rows (header + 5) remain.
$ws.Range("A7:A21").EntireRow.Delete()
